$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.058.62"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.621.35"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.02%  "
$savedStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.94"
$ws.Range("D5").Style = $savedStyle
$ws.Range("E5").Value = "  -1.33%  "
$savedStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("D6").Style = $savedStyle
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  -0.02%  "
$savedStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0627"
$ws.Range("D8").Style = $savedStyle
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  -1.61%  "
$savedStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.83"
$ws.Range("D10").Style = $savedStyle
$ws.Range("E10").Value = "  -0.70%  "
$savedStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = $savedStyle
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.851.38"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "1.624.81"
$ws.Range("E13").Value = "  -0.90%  "
$savedStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("D14").Style = $savedStyle
$ws.Range("E14").Value = "  -0.22%  "
$savedStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.538"
$ws.Range("D15").Style = $savedStyle
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "27.043.05"
$ws.Range("E16").Value = "  -0.25%  "
$savedStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.36"
$ws.Range("D17").Style = $savedStyle
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("E18").Value = "  -0.56%  "
$savedStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.66"
$ws.Range("D19").Style = $savedStyle
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("E20").Value = "  -0.10%  "
$savedStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("D21").Style = $savedStyle
$ws.Range("E21").Value = "  -1.29%  "
$savedStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("D22").Style = $savedStyle
$ws.Range("E22").Value = "  -2.04%  "
$savedStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("D23").Style = $savedStyle
$ws.Range("E23").Value = "  -7.43%  "
$savedStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.00"
$ws.Range("D24").Style = $savedStyle
$ws.Range("E24").Value = "  -1.05%  "
$savedStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.06"
$ws.Range("D25").Style = $savedStyle
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  -0.06%  "
$savedStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.33"
$ws.Range("D27").Style = $savedStyle
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("E28").Value = "  -2.91%  "
$savedStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("D29").Style = $savedStyle
$ws.Range("E29").Value = "  -1.27%  "
$savedStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("D30").Style = $savedStyle
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  -0.82%  "
$savedStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").Style = $savedStyle
$ws.Range("E32").Value = "  -1.74%  "
$savedStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.705"
$ws.Range("D33").Style = $savedStyle
$ws.Range("E33").Value = "  +29.98%  "
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").Value = "1.351.51"
$ws.Range("E35").Value = "  +3.36%  "
$savedStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("D36").Style = $savedStyle
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("E38").Value = "  +0.05%  "
$savedStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.839"
$ws.Range("D39").Style = $savedStyle
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  +0.45%  "
$savedStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.796"
$ws.Range("D42").Style = $savedStyle
$ws.Range("E42").Value = "  -1.90%  "
$savedStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.33"
$ws.Range("D43").Style = $savedStyle
$ws.Range("E43").Value = "  +0.68%  "
$savedStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.79"
$ws.Range("D44").Style = $savedStyle
$ws.Range("E44").Value = "  +3.54%  "
$ws.Range("D45").Value = "1.761.53"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$savedStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.81"
$ws.Range("D46").Style = $savedStyle
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$savedStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.64"
$ws.Range("D47").Style = $savedStyle
$ws.Range("E47").Value = "  +2.76%  "
$savedStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.852"
$ws.Range("D48").Style = $savedStyle
$ws.Range("E48").Value = "  +27.69%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  +4.47%  "
$ws.Range("E51").Value = "  +0.16%  "
